$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2024-06-02 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-03 Monday", 2) | Out-Null

# Update the division problems/answers in the table, cell by cell
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "765÷5=153, 0"
$t.Cell(1, 2).Range.Text = "605÷7=86, 3"
$t.Cell(1, 3).Range.Text = "893÷8=111, 5"
$t.Cell(1, 4).Range.Text = "530÷8=66, 2"
$t.Cell(1, 5).Range.Text = "602÷8=75, 2"
$t.Cell(5, 1).Range.Text = "110÷9=12, 2"
$t.Cell(5, 2).Range.Text = "888÷2=444, 0"
$t.Cell(5, 3).Range.Text = "285÷9=31, 6"
$t.Cell(5, 4).Range.Text = "871÷2=435, 1"
$t.Cell(5, 5).Range.Text = "513÷7=73, 2"
$t.Cell(9, 1).Range.Text = "832÷4=208, 0"
$t.Cell(9, 2).Range.Text = "793÷3=264, 1"
$t.Cell(9, 3).Range.Text = "303÷4=75, 3"
$t.Cell(9, 4).Range.Text = "990÷8=123, 6"
$t.Cell(9, 5).Range.Text = "117÷2=58, 1"
$t.Cell(13, 1).Range.Text = "769÷3=256, 1"
$t.Cell(13, 2).Range.Text = "643÷9=71, 4"
$t.Cell(13, 3).Range.Text = "991÷6=165, 1"
$t.Cell(13, 4).Range.Text = "613÷7=87, 4"
$t.Cell(13, 5).Range.Text = "100÷7=14, 2"
$t.Cell(17, 1).Range.Text = "191÷5=38, 1"
$t.Cell(17, 2).Range.Text = "370÷4=92, 2"
$t.Cell(17, 3).Range.Text = "869÷4=217, 1"
$t.Cell(17, 4).Range.Text = "995÷8=124, 3"
$t.Cell(17, 5).Range.Text = "567÷8=70, 7"
